$wb = $excel.ActiveWorkbook

# --- Rename the "Include from Vaccine Gender" sheet to "Include #0" ---
$includeSheet = $wb.Worksheets.Item("Include from Vaccine Gender")
$includeSheet.Name = "Include #0"

# --- Metadata sheet updates ---
$meta = $wb.Worksheets.Item("Metadata")

# Update the ValueSet URL (pythia -> cicada)
$meta.Range("B2").Value = "http://fhirfli.dev/fhir/ig/cicada/ValueSet/vaccine-gender"

# Update the Date value
$meta.Range("B8").Value = "2026-02-11T14:37:07-05:00"

# Insert a new "Jurisdiction" row after "Contact" (row 10), before "Description" (row 11)
$meta.Rows.Item(11).Insert()

# Match the formatting used by the other data rows (the insert leaves the new
# row with default formatting, so pull it from the row right below).
$meta.Range("A12:B12").Copy()
$meta.Range("A11:B11").PasteSpecial(-4122)

$meta.Range("A11").Value = "Jurisdiction"
$meta.Range("B11").Value = $includeSheet.Range("A6").Value

# --- Include sheet updates ---
$includeSheet.Range("B7").Value = "http://fhirfli.dev/fhir/ig/cicada/CodeSystem/VaccineGender"
